$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Slide 1's title placeholder ("Title 1", shape id 2) currently reads
# "Name of This Business" - replace it with the new business name.
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "So-Show"
